$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 101 (shifts old rows 101..164 down to 102..165)
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new record
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44582
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100103
$ws.Range("H101").Value = "Frutos de hueso (carozo)"
$ws.Range("I101").Value = 100103002
$ws.Range("J101").Value = "Ciruela"
$ws.Range("K101").Value = "Black Amber"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 85
$ws.Range("N101").Value = 12000
$ws.Range("O101").Value = 13000
$ws.Range("P101").Value = 12412
$ws.Range("Q101").Value = '$/bandeja 18 kilos granel'
$ws.Range("R101").Value = "Región de O'Higgins"
$ws.Range("S101").Value = 690
$ws.Range("T101").Value = 18
